# Tech-report title page: collapse one redundant blank "Title page: text"
# paragraph out of each blank run (between the number/year/TITLE/Authors/
# Address bookmarks) and left-justify the paragraph that absorbs its place.
#
# Each of the four blank runs starts life as 6 empty "Titlepagetext"
# paragraphs; after the edit it is 5, and the paragraph immediately
# preceding the next labeled line carries an explicit <w:jc w:val="left"/>
# (the style's own default justification is "center").

$d = $word.ActiveDocument

# Anchor on the bookmarks rather than raw paragraph numbers so the script
# keys off the document's actual content/layout.
$anchors = @("number", "year", "title", "authors", "address")
$starts = @()
foreach ($name in $anchors) {
    $bm = $d.Bookmarks.Item($name)
    $starts += $bm.Range.Paragraphs.Item(1).Index
}

# (deleteOffset, alignOffset) from each run's starting bookmark paragraph -
# identifies exactly which pair of the run's identical blank paragraphs
# gets merged, matching the source edit paragraph-for-paragraph.
$offsets = @(
    @(4, 5),
    @(2, 3),
    @(3, 4),
    @(2, 3)
)

# Apply from the last run to the first so deleting earlier paragraphs in
# later runs never invalidates the paragraph indices still to be visited.
for ($k = $offsets.Count - 1; $k -ge 0; $k--) {
    $runStart = $starts[$k]
    $deleteIndex = $runStart + $offsets[$k][0]
    $alignIndex = $runStart + $offsets[$k][1]

    $d.Paragraphs.Item($alignIndex).Format.Alignment = 0  ; # wdAlignParagraphLeft
    $d.Paragraphs.Item($deleteIndex).Range.Delete()
}
